# Update the two time-slot values in column C (rows 2 and 3).
# Column B (the encrypted/coded labels) and the rest of the table
# (rows 4-7) stay exactly as they were.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "8:20-8:25"
$ws.Range("C3").Value = "8:25-8:30"

# Move the active cell/selection from B11 down to B12, matching the
# sheet's saved cursor position.
$ws.Range("B12").Select()
